$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data appended to the bottom of the sheet (intentional overcontrol fault rows)
$data = @(
    @(100, -2.88333003678403),
    @(101, 2.53389756171697),
    @(102, -2.45399065399328),
    @(103, -7.85260295716882),
    @(104, 2.44345098578806),
    @(105, -1.24246538667722),
    @(106, 2.03668545972892),
    @(107, -0.363658423954536),
    @(108, 3.16646918371305),
    @(109, -2.61257155538862),
    @(110, 2.41028056501556),
    @(111, -2.88333003678403),
    @(112, 2.53389756171697),
    @(113, -2.45399065399328),
    @(114, 2.44345098578806),
    @(115, -1.24246538667722),
    @(116, 2.03668545972892),
    @(117, -0.363658423954536)
)

$startRow = 105
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Update view / selection state to match the scrolled-down position after
# appending the new rows.
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.ScrollRow = 88
$ws.Range("A123").Select()
